# This script reproduces a small, surgical set of run-structure edits that
# mirror what a human author produced when lightly re-touching two bullet
# lines and a date range in this resume (the XML-level diff shows runs being
# re-split/merged around the point the cursor last visited -- i.e. Word's
# auto "_GoBack" bookmark moving -- plus the removal of a stale grammar
# proofing mark around "October 2013").
#
# Helper pattern used below: to force Word to re-coalesce two (or more)
# adjacent, identically-formatted runs into a single run, we briefly place
# "fence" bookmarks just outside the span we want merged (so the
# re-coalescing can't spill past them), make a throw-away text edit inside
# the span (Word only recombines runs where text actually changed), put the
# original text back, then remove the fence bookmarks.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "using Shi" + "ny &"  ->  single run "using Shiny &"
#    (this also discards the "_GoBack" bookmark that used to sit between
#    them, since it's interior to the merged span)
# ---------------------------------------------------------------------
$fenceL = $d.Content
$fenceL.Find.Execute("using Shi", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fenceL.Collapse(1)
$d.Bookmarks.Add("ZZFenceL", $fenceL) | Out-Null

$fenceR = $d.Content
$fenceR.Find.Execute("using Shiny &", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fenceR.Collapse(0)
$d.Bookmarks.Add("ZZFenceR", $fenceR) | Out-Null

$work = $d.Content
$work.Find.Execute("using Shiny &", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$work.Text = "using ShinyZZZTMP &"

$work2 = $d.Content
$work2.Find.Execute("using ShinyZZZTMP &", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$work2.Text = "using Shiny &"

$d.Bookmarks("ZZFenceL").Delete()
$d.Bookmarks("ZZFenceR").Delete()

# ---------------------------------------------------------------------
# 2) "October" + " 2013"  ->  single run "October 2013"
#    (this also removes the now-orphaned gramStart/gramEnd proofing marks
#    that bracketed "October")
# ---------------------------------------------------------------------
$fenceL2 = $d.Content
$fenceL2.Find.Execute("October", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fenceL2.Collapse(1)
$d.Bookmarks.Add("ZZFenceL2", $fenceL2) | Out-Null

$fenceR2 = $d.Content
$fenceR2.Find.Execute("October 2013", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fenceR2.Collapse(0)
$d.Bookmarks.Add("ZZFenceR2", $fenceR2) | Out-Null

$work3 = $d.Content
$work3.Find.Execute("October 2013", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$work3.Text = "OctoberZZZTMP 2013"

$work4 = $d.Content
$work4.Find.Execute("OctoberZZZTMP 2013", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$work4.Text = "October 2013"

$d.Bookmarks("ZZFenceL2").Delete()
$d.Bookmarks("ZZFenceR2").Delete()

# ---------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark (Word's "last edit location" marker) to
#    sit inside "Create data and reporting standards, shift " -- right
#    after "Create data a" -- splitting that run in two, the way Word
#    would leave it after the cursor's final edit there.
# ---------------------------------------------------------------------
$goback = $d.Content
$goback.Find.Execute("Create data a", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goback.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goback) | Out-Null
